$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(4, 6, 7, 8)
foreach ($r in $rows) {
    $ws.Range("C$r").Value = "-"
    $ws.Range("D$r").Value = "MEC-2B-Metalografia"
}
